# Apply the ORM_Records.xlsx changes:
#  - Add a new "Remark" column (G) to the table, with a couple of column-width tweaks
#  - Record the Closed Date for the first entry (QualityPlan.docx)
#  - Add two new ORM records for "TUCMS.docs" and "HLD.docx" reviews (rows 11-12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data rows (No=10 / TUCMS.docs, No=11 / HLD.docx) -----------------
# Row 11: ORM No 31081810, TUCMS.docs, Kaung Myat Bo, issued 8/31/2018, remark "After 2nd Audit"
$ws.Range("B11").Value = 31081810
$ws.Range("C11").Value = "TUCMS.docs"
$ws.Range("D11").Value = "Kaung Myat Bo"
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)   # xlPasteFormats (reuse existing date style)
$ws.Range("E11").Value = 43343
$ws.Range("G11").Value = "After 2nd Audit"

# Row 12: ORM No 31081811, HLD.docx, Kaung Myat Bo, issued 8/31/2018
$ws.Range("B12").Value = 31081811
$ws.Range("C12").Value = "HLD.docx"
$ws.Range("D12").Value = "Kaung Myat Bo"
$ws.Range("E2").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = 43343

# --- Closed Date for the first record (QualityPlan.docx) ------------------
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = 43303
$ws.Range("E2").Copy()
$ws.Range("G2").PasteSpecial(-4122)    # keep the same date style, but leave it blank

# --- Extend the table with a new "Remark" column ---------------------------
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()
$ws.Range("F1").Copy()
$col.Range.Cells.Item(1, 1).PasteSpecial(-4122)   # reuse the Accent1 header style
$col.Range.Cells.Item(1, 1).Value = "Remark"

# Match formatting/border style already used at the bottom of the table
$ws.Range("F16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("G23").PasteSpecial(-4122)

# --- Column widths -----------------------------------------------------------
$fWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(7).ColumnWidth = $fWidth
$ws.Columns.Item(8).ColumnWidth = 14.6

# --- Selection shown in the saved view --------------------------------------
$ws.Range("I4").Select()
